$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 8: IBL_Irradiance
$ws.Range("A8").Value = "IBL_Irradiance"
$ws.Range("B8").Value = 1
$ws.Range("C8").Value = "WrapLinear"
$ws.Range("D8").Value = "Sampler"
$ws.Range("F8").Value = "`"PS`""
$ws.Range("G8").Value = 0
$ws.Range("H8").Value = 1

# Row 9: IBL_Radiance
$ws.Range("A9").Value = "IBL_Radiance"
$ws.Range("B9").Value = 1
$ws.Range("C9").Value = "WrapLinear"
$ws.Range("D9").Value = "Sampler"
$ws.Range("F9").Value = "`"PS`""
$ws.Range("G9").Value = 0
$ws.Range("H9").Value = 1

# Update selection to match the diff (active cell C7)
$ws.Range("C7").Select()
